$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "86.177.33"
Set-TextValue $ws.Range("E2") "  +8.14%  "

Set-TextValue $ws.Range("D3") "3.302.95"
Set-TextValue $ws.Range("E3") "  +4.33%  "

Set-TextValue $ws.Range("E4") "  +0.16%  "

Set-TextValue $ws.Range("D5") "217.76"
Set-TextValue $ws.Range("E5") "  +4.90%  "

Set-TextValue $ws.Range("D6") "635.81"
Set-TextValue $ws.Range("E6") "  +1.19%  "

Set-TextValue $ws.Range("D7") "0.318"
Set-TextValue $ws.Range("E7") "  +18.00%  "

Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.02%  "

Set-TextValue $ws.Range("D9") "0.596"
Set-TextValue $ws.Range("E9") "  +0.10%  "

Set-TextValue $ws.Range("D10") "3.305.54"
Set-TextValue $ws.Range("E10") "  +4.33%  "

Set-TextValue $ws.Range("D11") "0.599"
Set-TextValue $ws.Range("E11") "  -1.39%  "

Set-TextValue $ws.Range("E12") "  +6.98%  "

Set-TextValue $ws.Range("E13") "  +0.89%  "

Set-TextValue $ws.Range("D14") "3.914.05"
Set-TextValue $ws.Range("E14") "  +4.12%  "

Set-TextValue $ws.Range("D15") "34.02"
Set-TextValue $ws.Range("E15") "  +6.42%  "

Set-TextValue $ws.Range("D16") "5.36"
Set-TextValue $ws.Range("E16") "  +1.41%  "

Set-TextValue $ws.Range("D17") "85.739.50"

Set-TextValue $ws.Range("D18") "3.307.77"
Set-TextValue $ws.Range("E18") "  +3.65%  "

Set-TextValue $ws.Range("D19") "14.54"
Set-TextValue $ws.Range("E19") "  +1.08%  "

Set-TextValue $ws.Range("E20") "  +7.13%  "

Set-TextValue $ws.Range("D21") "444.47"
Set-TextValue $ws.Range("E21") "  +0.49%  "

Set-TextValue $ws.Range("D22") "9.10"

Set-TextValue $ws.Range("D23") "5.23"
Set-TextValue $ws.Range("E23") "  -1.60%  "

Set-TextValue $ws.Range("D24") "7.40"
Set-TextValue $ws.Range("E24") "  +7.33%  "

Set-TextValue $ws.Range("D25") "5.37"
Set-TextValue $ws.Range("E25") "  +13.76%  "

Set-TextValue $ws.Range("D26") "12.20"
Set-TextValue $ws.Range("E26") "  +11.82%  "

Set-TextValue $ws.Range("D27") "3.489.10"
Set-TextValue $ws.Range("E27") "  +4.29%  "

Set-TextValue $ws.Range("D28") "78.17"
Set-TextValue $ws.Range("E28") "  +1.84%  "

Set-TextValue $ws.Range("E29") "  +5.78%  "

Set-TextValue $ws.Range("E30") "  -0.04%  "

Set-TextValue $ws.Range("D31") "0.174"
Set-TextValue $ws.Range("E31") "  +42.66%  "

Set-TextValue $ws.Range("D32") "602.50"
Set-TextValue $ws.Range("E32") "  +9.71%  "

Set-TextValue $ws.Range("E33") "  +0.57%  "

Set-TextValue $ws.Range("D34") "1.00"
Set-TextValue $ws.Range("E34") "  -0.09%  "

Set-TextValue $ws.Range("E36") "  +1.96%  "

Set-TextValue $ws.Range("D37") "0.149"
Set-TextValue $ws.Range("E37") "  -0.38%  "

Set-TextValue $ws.Range("D38") "23.30"
Set-TextValue $ws.Range("E38") "  -0.10%  "

Set-TextValue $ws.Range("D39") "6.46"
Set-TextValue $ws.Range("E39") "  +13.88%  "

Set-TextValue $ws.Range("D40") "0.999"
Set-TextValue $ws.Range("E40") "  -0.18%  "

Set-TextValue $ws.Range("E41") "  +1.74%  "

Set-TextValue $ws.Range("D42") "21.30"
Set-TextValue $ws.Range("E42") "  +2.53%  "

Set-TextValue $ws.Range("D43") "3.12"
Set-TextValue $ws.Range("E43") "  +15.68%  "

Set-TextValue $ws.Range("D44") "2.06"
Set-TextValue $ws.Range("E44") "  +13.11%  "

Set-TextValue $ws.Range("B45") "Monero"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D45") "158.64"
Set-TextValue $ws.Range("E45") "  -3.70%  "

Set-TextValue $ws.Range("B46") "USDe"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "1.00"
Set-TextValue $ws.Range("E46") "  +0.02%  "

Set-TextValue $ws.Range("D47") "187.64"
Set-TextValue $ws.Range("E47") "  -0.75%  "

Set-TextValue $ws.Range("E48") "  +3.66%  "

Set-TextValue $ws.Range("D49") "45.24"
Set-TextValue $ws.Range("E49") "  +3.95%  "

Set-TextValue $ws.Range("D50") "0.783"
Set-TextValue $ws.Range("E50") "  -0.28%  "

Set-TextValue $ws.Range("D51") "26.20"
Set-TextValue $ws.Range("E51") "  +4.68%  "
